$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 with new feature request
$ws.Range("A15").Value = "Minimum roll value"
$ws.Range("B15").Value = "When you roll below a threshold, treat any roll that is less than that threshold as that threshold."
$ws.Range("D15").Value = "Weston Fiala"

# Update the selected cell on the sheet
$ws.Range("B28").Select()
